$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Genre: Adventure, Puzzle -> Genre: Arcade
# ------------------------------------------------------------------
$d.Content.Find.Execute("Genre: Adventure, Puzzle", $true, $false, $false, $false, $false, $true, 1, $false, "Genre: Arcade", 2)

# ------------------------------------------------------------------
# 2. "The Adventure game based on movie ..." -> "The Arcade game based on movie ..."
# ------------------------------------------------------------------
$d.Content.Find.Execute("The Adventure game based on movie Stuart Little (1999) in style of 70-80s Atari 2600 games. ", $true, $false, $false, $false, $false, $true, 1, $false, "The Arcade game based on movie Stuart Little (1999) in style of 70-80s Atari 2600 games. ", 2)

# ------------------------------------------------------------------
# 3. Remove the whole paragraph describing the gameplay levels
#    ("You, in role of Stuart himself, ...") -- it disappears entirely,
#    merging into the (already present) following empty paragraph.
# ------------------------------------------------------------------
$gameplayText = "You, in role of Stuart himself, complete levels which were key moments from the film, Like meeting the Stuart, finding costumes, Escape from cats, Ship Race, Car Race and etc."
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Contains($gameplayText)) {
        $p.Range.Delete()
        break
    }
}

# ------------------------------------------------------------------
# 4. Target audience text
# ------------------------------------------------------------------
$d.Content.Find.Execute("*Target audience:  People who like classic adventures and puzzles (6+)", $true, $false, $false, $false, $false, $true, 1, $false, "*Target audience:  People who like classic Arcades and Atari style (6+)", 2)

# ------------------------------------------------------------------
# 5. Roadmap date fix: 20.04.24 -> 21.04.24
# ------------------------------------------------------------------
$d.Content.Find.Execute("20.04.24-polishing the game and sending it to jam.", $true, $false, $false, $false, $false, $true, 1, $false, "21.04.24-polishing the game and sending it to jam.", 2)

# ------------------------------------------------------------------
# 6. Move the "_GoBack" bookmark from the end of the document (after
#    "**Token Sabit - Tester") to right after the "Description: " run.
#    A zero-width Range sitting exactly on a paragraph mark confuses
#    Bookmarks.Add, so we briefly insert a placeholder character,
#    bookmark the 1-character Range around it, then delete that
#    character -- the (now empty) bookmark collapses into place and
#    Bookmarks.Add("_GoBack", ...) automatically relocates the
#    pre-existing bookmark of the same name.
# ------------------------------------------------------------------
$descRange = $d.Content
$descRange.Find.Execute("Description: ")
$descRange.Collapse(0)
$insertPos = $descRange.Start
$descRange.InsertBefore("X")
$placeholder = $d.Range($insertPos, $insertPos + 1)
$d.Bookmarks.Add("_GoBack", $placeholder)
$d.Bookmarks.Item("_GoBack").Range.Delete()
